$wb = $excel.ActiveWorkbook

# --- Sheet: Battery_Data ---
$ws1 = $wb.Worksheets.Item("Battery_Data")
$ws1.Range("B2").Value2 = 508909.30501200003
$ws1.Range("B3").Value2 = 376592.88570888003
$ws1.Range("B4").Value2 = 7531.8577141776004
$ws1.Range("B5").Value2 = 50418.5832974

# --- Sheet: Yearly BRC ---
$ws2 = $wb.Worksheets.Item("Yearly BRC")
$ws2.Range("B2").Value2 = 25202.447014969581
$ws2.Range("B3").Value2 = 25216.136282408941

# --- Rename "Reposition" -> "Replacement" across all sheets (keeps shared-string ordering stable) ---
foreach ($ws in $wb.Worksheets) {
    [void]$ws.Cells.Replace("Reposition", "Replacement")
}

# --- Re-fit column A widths to the (now longer) label text, mirroring Excel's bestFit behaviour ---
$ws1.Columns.Item(1).ColumnWidth = 36
$ws2.Columns.Item(1).ColumnWidth = 28.75
